# Update Name of Algo
# Applies updated RandomForest imputation values to the result_data worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 9.121399999999996
$ws.Range("B7").Value = 5.080200000000002
$ws.Range("E7").Value = 15.5173
$ws.Range("E15").Value = 16.2418
$ws.Range("B16").Value = 6.7318
$ws.Range("E21").Value = 16.86880000000001
$ws.Range("E22").Value = 16.8642
$ws.Range("E23").Value = 16.14569999999999
$ws.Range("B28").Value = 5.982600000000002
$ws.Range("B29").Value = 4.909000000000003
$ws.Range("B32").Value = 7.440699999999996
$ws.Range("E34").Value = 17.02800000000001
$ws.Range("B40").Value = 8.850299999999999
$ws.Range("E43").Value = 17.29440000000001
$ws.Range("E45").Value = 16.3945
$ws.Range("E50").Value = 16.3431
$ws.Range("E51").Value = 17.3183
$ws.Range("B52").Value = 5.165699999999999
$ws.Range("B57").Value = 4.890299999999997
$ws.Range("B66").Value = 5.648899999999998
$ws.Range("E66").Value = 17.08810000000001
$ws.Range("E67").Value = 17.1998
$ws.Range("E79").Value = 18.17260000000002
$ws.Range("E84").Value = 16.52659999999999
$ws.Range("E92").Value = 18.44770000000001
$ws.Range("E97").Value = 16.605
$ws.Range("B100").Value = 5.633999999999999
